# Adds an "exchange_type" column to the "Add Exchanges" sheet, filled with
# "technosphere" for every existing data row, and restores the selections
# that Excel records after that edit (as seen in the commit's XML diff).

$wb = $excel.ActiveWorkbook

# --- "Create Activities" sheet: cursor moved to H1 -------------------------
$ws1 = $wb.Worksheets.Item("Create Activities")
$ws1.Activate() | Out-Null
$ws1.Range("H1").Select() | Out-Null

# --- "Add Exchanges" sheet: insert the new exchange_type column -----------
$ws2 = $wb.Worksheets.Item("Add Exchanges")
$ws2.Activate() | Out-Null

# Inserting a whole column before column J (10) shifts the old J:N block to
# K:O, carrying cell styles/values along -- exactly matching the diff, which
# shows the previous J1 ("exchange_code") becoming K1, etc.
$ws2.Columns.Item(10).Insert() | Out-Null

# New header cell for the inserted column.
$ws2.Cells.Item(1, 10).Value = "exchange_type"

# Every data row (2-9) gets "technosphere" in the new column.
for ($r = 2; $r -le 9; $r++) {
    $ws2.Cells.Item($r, 10).Value = "technosphere"
}

# Final selection on this sheet, matching the diff.
$ws2.Range("F6").Select() | Out-Null
